$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Chad Bettis"
$ws.Range("B2").Value = "bettich01"
$ws.Range("C2").Value = "June 06 2017"
$ws.Range("D2").Value = "Illness"
$ws.Range("E2").Value = "Bettis is on the 60-day disabled list while recovering from testicular cancer but is expected to make his season debut sometime around the All-Star break."

$ws.Range("A3").Value = "David Dahl"
$ws.Range("B3").Value = "dahlda01"
$ws.Range("C3").Value = "June 23 2017"
$ws.Range("D3").Value = "Ribs"
$ws.Range("E3").Value = "Dahl is on the 10-day disabled list with a stress reaction of his sixth rib. It is unknown as to when he will be ready to join the lineup."

$ws.Range("A4").Value = "Carlos Gonzalez"
$ws.Range("B4").Value = "gonzaca01"
$ws.Range("C4").Value = "June 25 2017"
$ws.Range("D4").Value = "Shoulder"
$ws.Range("E4").Value = "Gonzalez has missed the last two games with a sore right shoulder and his status is uncertain for Sunday's game against the Dodgers."

$ws.Range("A5").Value = "Jon Gray"
$ws.Range("B5").Value = "grayjo02"
$ws.Range("C5").Value = "June 20 2017"
$ws.Range("D5").Value = "Toe"
$ws.Range("E5").Value = "Gray was placed on the 10-day disabled list with a stress fracture in his left foot. He is on a rehab assignment and is expected to rejoin the rotation before the beginning of July."

$ws.Range("A6").Value = "Gerardo Parra"
$ws.Range("B6").Value = "parrage01"
$ws.Range("C6").Value = "June 20 2017"
$ws.Range("D6").Value = "Quadricep"
$ws.Range("E6").Value = "Parra has landed on the 10-day disabled list with a strained right quadriceps and is likely to be sidelined until the start of July."

$ws.Range("A7").Value = "Chad Qualls"
$ws.Range("B7").Value = "quallch01"
$ws.Range("C7").Value = "June 21 2017"
$ws.Range("D7").Value = "Back"
$ws.Range("E7").Value = "Qualls has been placed on the 10-day disabled list with lower back spasms and it is unclear as to if he will return for Monday's game against the Giants."

$ws.Range("B20").Select()
